$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A317").Value = 316
$ws.Range("B317").Value = "Which of the following is the primary focus of Lean Portfolio Management?"
$ws.Range("C317").Value = "Define the enterprise strategy; Establish lean budgets; Align strategy and execution; Ensure program and portfolio execution"
$ws.Range("D317").Value = "Align strategy and execution"
$ws.Range("E317").Value = 0
$ws.Range("F317").Value = 0
$ws.Range("A318").Value = 317
$ws.Range("B318").Value = "What is the primary purpose of Strategic Themes?"
$ws.Range("C318").Value = "Determine the order in which Epics should be executed; Drive incremental implementation across the enterprise; Define the sequence of steps used to deliver value to the customer; Connect the portfolio to the enterprise business strategy"
$ws.Range("D318").Value = "Connect the portfolio to the enterprise business strategy"
$ws.Range("E318").Value = 0
$ws.Range("F318").Value = 0
$ws.Range("A319").Value = 318
$ws.Range("B319").Value = "What is one fundamental difference between Agile Software Engineering and Waterfall?"
$ws.Range("C319").Value = "Agile delivers every single requested requirement.; Waterfall has no fixed schedule.; Agile supports test-first development.; Waterfall concentrates on code quality."
$ws.Range("D319").Value = "Agile supports test-first development."
$ws.Range("E319").Value = 0
$ws.Range("F319").Value = 0
$ws.Range("A320").Value = 319
$ws.Range("B320").Value = "What is the purpose of having deployment decoupled from release?"
$ws.Range("C320").Value = "To enable deploying upon demand.; To remove the need to respond quickly to product issues.; To enable releasing upon demand.; To make deploying a business decision."
$ws.Range("D320").Value = "To enable releasing upon demand."
$ws.Range("E320").Value = 0
$ws.Range("F320").Value = 0
$ws.Range("A321").Value = 320
$ws.Range("B321").Value = "Which of the following is a Lean Budget Guardrail?"
$ws.Range("C321").Value = "Participatory budgeting; Investment horizons; Centralized decision-making; Lean startup cycle"
$ws.Range("D321").Value = "Investment horizons"
$ws.Range("E321").Value = 0
$ws.Range("F321").Value = 0
$ws.Range("A322").Value = 321
$ws.Range("B322").Value = "What is BDD?"
$ws.Range("C322").Value = "Business-Driven Development; Benchmark-Driven Development; Behavior-Driven Development; Backlog-Driven Development"
$ws.Range("D322").Value = "Behavior-Driven Development"
$ws.Range("E322").Value = 0
$ws.Range("F322").Value = 0
$ws.Range("A323").Value = 322
$ws.Range("B323").Value = "Which two actions can slow the delivery of value?"
$ws.Range("C323").Value = "Shortened Architectural Runway; Loopbacks; System delays; Team swarming"
$ws.Range("D323").Value = "Loopbacks; System delays"
$ws.Range("E323").Value = 0
$ws.Range("F323").Value = 0
$ws.Range("A324").Value = 323
$ws.Range("B324").Value = "Which two statements best describe a cross-functional Agile Team?"
$ws.Range("C324").Value = "They are optimized for communication and delivery of value.; They can define, build, and test a Feature or component.; They deliver value every six weeks.; They release customer products to production continuously.; They are made up of members, each of whom can define, develop, test, and deploy the system."
$ws.Range("D324").Value = "They are optimized for communication and delivery of value.; They can define, build, and test a Feature or component."
$ws.Range("E324").Value = 0
$ws.Range("F324").Value = 0
$ws.Range("A325").Value = 324
$ws.Range("B325").Value = "The role of the Product Manager is most similar to what other role in SAFe?"
$ws.Range("C325").Value = "Development Manager; Solution Manager; Release Train Engineer; Business Manager"
$ws.Range("D325").Value = "Solution Manager"
$ws.Range("E325").Value = 0
$ws.Range("F325").Value = 0
$ws.Range("A326").Value = 325
$ws.Range("B326").Value = "What is the foundation of Lean?"
$ws.Range("C326").Value = "Innovation; Lean-Agile development; Leadership; Lean-thinking"
$ws.Range("D326").Value = "Leadership"
$ws.Range("E326").Value = 0
$ws.Range("F326").Value = 0
$ws.Range("A327").Value = 326
$ws.Range("B327").Value = "What is a benefit of continuously deploying using a DevOps model?"
$ws.Range("C327").Value = "It alleviates the reliance on the skill sets of Agile Teams.; It ensures that changes deployed to production are always immediately available to end-users.; It lessens the severity and frequency of release failures.; It increases the transaction cost."
$ws.Range("D327").Value = "It lessens the severity and frequency of release failures."
$ws.Range("E327").Value = 0
$ws.Range("F327").Value = 0
$ws.Range("A328").Value = 327
$ws.Range("B328").Value = "What are three practices for building large Solutions?"
$ws.Range("C328").Value = "Employ and improve the Continuous Development Pipeline.; Employ multiple Solution Teams.; Build Solution components and capabilities with ARTs.; Build and integrate the Solution with a Solution Train.; Ensure value stream consistency.; Capture and refine system specifications as fixed/variable Solution Intent."
$ws.Range("D328").Value = "Build and integrate the Solution with a Solution Train.; Build Solution components and capabilities with ARTs.; Capture and refine system specifications as fixed/variable Solution Intent."
$ws.Range("E328").Value = 0
$ws.Range("F328").Value = 0
$ws.Range("A329").Value = 328
$ws.Range("B329").Value = "What is the ultimate goal of DevOps in SAFe?"
$ws.Range("C329").Value = "To enable Agile Release Trains to deliver value more frequently.; To combine Development and Operations groups into one team.; To clarify program and team roles.; To automate the integration, testing, and deployment of features."
$ws.Range("D329").Value = "To enable Agile Release Trains to deliver value more frequently."
$ws.Range("E329").Value = 0
$ws.Range("F329").Value = 0
$ws.Range("A330").Value = 329
$ws.Range("B330").Value = "In which event are the PI Objectives created?"
$ws.Range("C330").Value = "Post PI-Planning; PI Planning; Iteration Planning; PI System Demo"
$ws.Range("D330").Value = "PI Planning"
$ws.Range("E330").Value = 0
$ws.Range("F330").Value = 0
$ws.Range("A331").Value = 330
$ws.Range("B331").Value = "Which statement is most accurate about the Solution Vision?"
$ws.Range("C331").Value = "It summarizes the team PI Objectives for the current Program Increment.; It provides an outline of the Features for the next three Program Increments.; It expresses the strategic intent of the Program.; It drives the allocation of budget for the Agile Release Train."
$ws.Range("D331").Value = "It expresses the strategic intent of the Program."
$ws.Range("E331").Value = 0
$ws.Range("F331").Value = 0
$ws.Range("A332").Value = 331
$ws.Range("B332").Value = "What does Solution Intent describe?"
$ws.Range("C332").Value = "The customer or consumer of the solution.; The behavior of the solution.; The platform architecture.; Where the solution operates."
$ws.Range("D332").Value = "The behavior of the solution."
$ws.Range("E332").Value = 0
$ws.Range("F332").Value = 0
$ws.Range("A333").Value = 332
$ws.Range("B333").Value = "Which of the following is a key purpose of DevOps?"
$ws.Range("C333").Value = "DevOps joins development and operations to enable continuous delivery.; DevOps focuses on automating the delivery pipeline to reduce transaction cost.; DevOps enables continuous release by building a scalable Continuous Delivery Pipeline.; DevOps focuses on a set of practices applied to large systems."
$ws.Range("D333").Value = "DevOps joins development and operations to enable continuous delivery."
$ws.Range("E333").Value = 0
$ws.Range("F333").Value = 0
$ws.Range("A334").Value = 333
$ws.Range("B334").Value = "What is the purpose of dynamic models?"
$ws.Range("C334").Value = "They are used when there are numerous and complex interactions.; They are used when there are a set of parallel interactions related to a scenario.; They are used when there are a specific number of interactions that are less complex.; They are used to sequence a set of interactions related to a scenario."
$ws.Range("D334").Value = "They are used when there are numerous and complex interactions."
$ws.Range("E334").Value = 0
$ws.Range("F334").Value = 0
$ws.Range("A335").Value = 334
$ws.Range("B335").Value = "At what level of the SAFe Big Picture do Strategic Themes reside?"
$ws.Range("C335").Value = "Essential Level.; Team Level.; Portfolio Level.; Large Solution Level."
$ws.Range("D335").Value = "Portfolio Level"
$ws.Range("E335").Value = 0
$ws.Range("F335").Value = 0
$ws.Range("A336").Value = 335
$ws.Range("B336").Value = "Which of the following is responsible for managing the Portfolio Kanban?"
$ws.Range("C336").Value = "Lean Portfolio Management.; Release Train Engineer.; Product Management.; Solution Management."
$ws.Range("D336").Value = "Lean Portfolio Management"
$ws.Range("E336").Value = 0
$ws.Range("F336").Value = 0
$ws.Range("C287").Select()
